$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Faturamento diario - dados atualizados para incluir maio/2025 (ADD 06-05-2025)
# e remover o mes de janeiro/2025, que saiu da janela de 4 meses exibida.

$data = @(
  @(2, 6309.66, 5, 2025, "05/2025"),
  @(5, 6038.94, 5, 2025, "05/2025"),
  @(6, 1805.49, 5, 2025, "05/2025"),
  @(1, 7890.54, 4, 2025, "04/2025"),
  @(2, 7622.33, 4, 2025, "04/2025"),
  @(3, 17544.87, 4, 2025, "04/2025"),
  @(4, 1905.42, 4, 2025, "04/2025"),
  @(7, 1044.7, 4, 2025, "04/2025"),
  @(8, 2415.17, 4, 2025, "04/2025"),
  @(9, 5303.6, 4, 2025, "04/2025"),
  @(10, 5184.47, 4, 2025, "04/2025"),
  @(11, 1553.36, 4, 2025, "04/2025"),
  @(14, 6895.1, 4, 2025, "04/2025"),
  @(15, 3115.18, 4, 2025, "04/2025"),
  @(16, 2957.12, 4, 2025, "04/2025"),
  @(17, 3488.72, 4, 2025, "04/2025"),
  @(22, 6143.24, 4, 2025, "04/2025"),
  @(23, 7854.64, 4, 2025, "04/2025"),
  @(24, 2184.55, 4, 2025, "04/2025"),
  @(25, 1904.36, 4, 2025, "04/2025"),
  @(28, 5438.35, 4, 2025, "04/2025"),
  @(29, 7952.53, 4, 2025, "04/2025"),
  @(30, 8340.92, 4, 2025, "04/2025"),
  @(3, 3051.74, 3, 2025, "03/2025"),
  @(6, 8076.58, 3, 2025, "03/2025"),
  @(7, 2846.82, 3, 2025, "03/2025"),
  @(10, 5273, 3, 2025, "03/2025"),
  @(11, 2859.88, 3, 2025, "03/2025"),
  @(12, 3586.54, 3, 2025, "03/2025"),
  @(13, 3667.31, 3, 2025, "03/2025"),
  @(14, 3928.66, 3, 2025, "03/2025"),
  @(17, 3320.22, 3, 2025, "03/2025"),
  @(18, 923.8099999999999, 3, 2025, "03/2025"),
  @(19, 4986.04, 3, 2025, "03/2025"),
  @(20, 1086.5, 3, 2025, "03/2025"),
  @(21, 6555.309999999999, 3, 2025, "03/2025"),
  @(24, 5885.61, 3, 2025, "03/2025"),
  @(25, 3667.99, 3, 2025, "03/2025"),
  @(26, 6259.42, 3, 2025, "03/2025"),
  @(27, 3653.25, 3, 2025, "03/2025"),
  @(28, 2642.05, 3, 2025, "03/2025"),
  @(31, 11719.82, 3, 2025, "03/2025"),
  @(3, 5672.85, 2, 2025, "02/2025"),
  @(4, 4327.82, 2, 2025, "02/2025"),
  @(5, 4515.13, 2, 2025, "02/2025"),
  @(6, 4392.82, 2, 2025, "02/2025"),
  @(7, 6569.1, 2, 2025, "02/2025"),
  @(10, 4368.17, 2, 2025, "02/2025"),
  @(11, 2353.68, 2, 2025, "02/2025"),
  @(12, 5981.65, 2, 2025, "02/2025"),
  @(13, 4046.06, 2, 2025, "02/2025"),
  @(14, 2145.6, 2, 2025, "02/2025"),
  @(17, 2525.47, 2, 2025, "02/2025"),
  @(18, 4856.68, 2, 2025, "02/2025"),
  @(19, 3002.33, 2, 2025, "02/2025"),
  @(20, 2102.41, 2, 2025, "02/2025"),
  @(21, 3130.24, 2, 2025, "02/2025"),
  @(24, 4723.690000000001, 2, 2025, "02/2025"),
  @(25, 3334.44, 2, 2025, "02/2025"),
  @(26, 18545.47, 2, 2025, "02/2025"),
  @(27, 10080.83, 2, 2025, "02/2025"),
  @(28, 13958.13, 2, 2025, "02/2025")
)

# Limpa todo o conteudo de dados existente (mantem o cabecalho na linha 1)
# para que a tabela de strings compartilhadas seja reconstruida do zero.
$ws.Range("A2:E81").ClearContents() | Out-Null

$r = 2
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
  $r++
}
